$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 77 (shifts existing rows 77.. down to 79..)
$ws.Rows.Item(77).Resize(2).Insert()

# Copy the date-number-format style from D79 (the row pushed down, which
# keeps its original style) onto the freshly inserted D77:D78 cells only.
$ws.Range("D79").Copy()
$ws.Range("D77:D78").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 77 data
$ws.Range("A77").Value = 7
$ws.Range("B77").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C77").Value = "Ñuble"
$ws.Range("D77").Value = 44452
$ws.Range("E77").Value = 16
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100101
$ws.Range("H77").Value = "Berries"
$ws.Range("I77").Value = 100101007
$ws.Range("J77").Value = "Kiwi"
$ws.Range("K77").Value = "Hayward"
$ws.Range("L77").Value = "Primera"
$ws.Range("M77").Value = 100
$ws.Range("N77").Value = 12000
$ws.Range("O77").Value = 12500
$ws.Range("P77").Value = 12250
$ws.Range("Q77").Value = "`$/bandeja 18 kilos"
$ws.Range("R77").Value = "Provincia de Curicó"
$ws.Range("S77").Value = 681
$ws.Range("T77").Value = 18

# Row 78 data
$ws.Range("A78").Value = 7
$ws.Range("B78").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C78").Value = "Ñuble"
$ws.Range("D78").Value = 44452
$ws.Range("E78").Value = 16
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100101
$ws.Range("H78").Value = "Berries"
$ws.Range("I78").Value = 100101007
$ws.Range("J78").Value = "Kiwi"
$ws.Range("K78").Value = "Hayward"
$ws.Range("L78").Value = "Segunda"
$ws.Range("M78").Value = 60
$ws.Range("N78").Value = 11000
$ws.Range("O78").Value = 11500
$ws.Range("P78").Value = 11250
$ws.Range("Q78").Value = "`$/bandeja 18 kilos"
$ws.Range("R78").Value = "Provincia de Curicó"
$ws.Range("S78").Value = 625
$ws.Range("T78").Value = 18

Write-Host "done"
